# Generate Report for Handback
#
# The handback transform failed for the "9b5efedf-10da-461f-8bcd-fcdfac65bdde"
# file in both the zh-cn and de-de locales (the generated handback file name
# did not match the handoff file name). Update the status for that row from
# "Ready for handoff" to "Handback transform failed" everywhere it is
# reported (Overview summary + each locale sheet), and record the failure
# detail in the "Error Detail" column of the locale sheets, widening that
# column so the message is readable.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: row 3 is the 9b5efedf-10da-461f-8bcd-fcdfac65bdde.md entry;
# column E is the zh-cn status, column F is the de-de status.
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Per-locale sheets: row 3 (Status column C) for the same file.
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Widen the "Error Detail" column (P) on both locale sheets so the failure
# message is readable.
$wsZhCn.Range("P1").ColumnWidth = 39.2
$wsDeDe.Range("P1").ColumnWidth = 39.2

# Record the handback/handoff file name mismatch in the Error Detail column.
$wsZhCn.Range("P3").Value = "Handback file name: o1b1sozd.s2n is different with handoff file name: 9b5efedf-10da-461f-8bcd-fcdfac65bdde.d57102120d156c86f19124bed81624c4da193120.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: o1b1sozd.s2n is different with handoff file name: 9b5efedf-10da-461f-8bcd-fcdfac65bdde.d57102120d156c86f19124bed81624c4da193120.de-de."
